$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# The original sheet mistakenly repeated the row-2 data values (company /
# policy name / owner) into the header row. Clear the whole used range and
# rewrite it from scratch: a proper header row, then the two data rows with
# their original company/policy/owner values plus the trailing metadata
# columns (property_category/category/date/legislator_name/legislator_id/
# source_file/index) that every other sheet in this workbook already has.
$ws.Range("A1:D3").ClearContents()

# --- Row 1: header ---
$ws.Cells.Item(1, 2).Value = "company"
$ws.Cells.Item(1, 3).Value = "name"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "property_category"
$ws.Cells.Item(1, 6).Value = "category"
$ws.Cells.Item(1, 7).Value = "date"
$ws.Cells.Item(1, 8).Value = "legislator_name"
$ws.Cells.Item(1, 9).Value = "legislator_id"
$ws.Cells.Item(1, 10).Value = "source_file"
$ws.Cells.Item(1, 11).Value = "index"

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = 103
$ws.Cells.Item(2, 2).Value = "富邦人壽"
$ws.Cells.Item(2, 3).Value = "安泰人壽靈活理財變額保險甲型"
$ws.Cells.Item(2, 4).Value = "賴士葆"
$ws.Cells.Item(2, 5).Value = "insurance"
$ws.Cells.Item(2, 6).Value = "normal"
$ws.Cells.Item(2, 7).Value = "2011-11-23"
$ws.Cells.Item(2, 8).Value = "賴士葆"
$ws.Cells.Item(2, 9).Value = 866
$ws.Cells.Item(2, 10).Value = "tmp2bc41"
$ws.Cells.Item(2, 11).Value = 103

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value = 104
$ws.Cells.Item(3, 2).Value = "國華人壽"
$ws.Cells.Item(3, 3).Value = "國華人壽終身壽險"
$ws.Cells.Item(3, 4).Value = "賴士葆"
$ws.Cells.Item(3, 5).Value = "insurance"
$ws.Cells.Item(3, 6).Value = "normal"
$ws.Cells.Item(3, 7).Value = "2011-11-23"
$ws.Cells.Item(3, 8).Value = "賴士葆"
$ws.Cells.Item(3, 9).Value = 866
$ws.Cells.Item(3, 10).Value = "tmp2bc41"
$ws.Cells.Item(3, 11).Value = 104

# Re-apply the header (bold + bordered) / data (bordered) styles that were
# cleared off column A, and extend them across the new E:K columns so the
# new cells match the look of the pre-existing B:D columns.
$ws.Range("A1").Style = $ws.Range("B1").Style
$ws.Range("A2:A3").Style = $ws.Range("B2").Style
$ws.Range("E1:K1").Style = $ws.Range("B1").Style
$ws.Range("E2:K3").Style = $ws.Range("B2").Style
